$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.017.38'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '2.404.09'
$ws.Range("E3").Value = '  -3.81%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '487.87'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.28'
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +19.94%  '
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.996'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '2.420.72'
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.30'
$ws.Range("E10").Value = '  +9.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0995'
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.335'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("E13").Value = '  +1.18%  '
$ws.Range("D14").Value = '2.826.16'
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '57.083.92'
$ws.Range("E15").Value = '  +0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.61'
$ws.Range("E16").Value = '  -3.20%  '
$ws.Range("E17").Value = '  -2.85%  '
$ws.Range("D18").Value = '2.415.05'
$ws.Range("E18").Value = '  -3.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").Value = '  +3.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '325.73'
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.96'
$ws.Range("E21").Value = '  -3.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.93'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.02'
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.405'
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.161'
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("D28").Value = '2.508.65'
$ws.Range("E28").Value = '  -3.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.24'
$ws.Range("E29").Value = '  -4.95%  '
$ws.Range("D30").Value = '0.0₃0781'
$ws.Range("E30").Value = '  -3.85%  '
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.90'
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.53'
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.27'
$ws.Range("E35").Value = '  -0.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.16'
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("E37").Value = '  -1.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.832'
$ws.Range("E38").Value = '  -4.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.102'
$ws.Range("E39").Value = '  +8.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '34.07'
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.52'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.995'
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.595'
$ws.Range("E44").Value = '  -3.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '269.30'
$ws.Range("E45").Value = '  -0.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0529'
$ws.Range("E46").Value = '  -6.10%  '
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("E48").Value = '  -0.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.51'
$ws.Range("E49").Value = '  -6.41%  '
$ws.Range("D50").Value = '1.871.60'
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.41'
$ws.Range("E51").Value = '  -3.20%  '
